# Apply updated price / volume(1h) figures for the crypto symbol list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "307.42"
    "E2" = "-0.24%"
    "D3" = "40.99"
    "E3" = "0.51%"
    "D4" = "5.232"
    "E4" = "2.26%"
    "D5" = "0.07661"
    "E5" = "0.61%"
    "D6" = "1.634"
    "E6" = "0.92%"
    "E7" = "1.61%"
    "D8" = "2.437"
    "E8" = "-0.45%"
    "D9" = "0.1245"
    "E9" = "13.28%"
    "D10" = "0.1824"
    "E10" = "3.06%"
    "D11" = "0.09091"
    "E11" = "-0.79%"
    "D12" = "0.04155"
    "E12" = "-0.32%"
    "D13" = "0.1049"
    "E13" = "-0.17%"
    "D14" = "0.001257"
    "E14" = "0.69%"
    "D15" = "0.005838"
    "E15" = "-0.76%"
    "E17" = "-0.25%"
    "E19" = "1.31%"
    "D20" = "7.513"
    "E20" = "13.77%"
    "D21" = "0.1383"
    "E21" = "1.37%"
    "D22" = "0.2882"
    "E22" = "7.45%"
    "D23" = "0.04068"
    "E23" = "0.34%"
    "E24" = "3.42%"
    "D25" = "0.004278"
    "E25" = "4.53%"
    "E26" = "-2.20%"
    "D38" = "0.02493"
    "E38" = "5.04%"
    "D39" = "0.05337"
    "E39" = "2.98%"
    "D40" = "0.007853"
    "E40" = "0.79%"
    "D41" = "0.1311"
    "E41" = "0.91%"
    "E42" = "1.94%"
    "E43" = "-1.95%"
    "D44" = "0.007666"
    "E44" = "-3.49%"
    "D45" = "0.3058"
    "E45" = "-0.58%"
    "D46" = "0.00006716"
    "E46" = "-3.23%"
    "D47" = "0.00000000751"
    "E47" = "0.04%"
    "D48" = "0.1699"
    "E48" = "442.06%"
    "D49" = "0.003105"
    "E49" = "-26.12%"
    "D50" = "0.00002103"
    "E50" = "0.04%"
    "D51" = "0.0002003"
    "E51" = "0.04%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text storage so values like "307.42" / "-0.24%" stay literal
    # strings (matching the source data) instead of being coerced to
    # numbers/percentages by Excel's automatic type detection.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Reset to the default (un-styled) cell style so we don't leave a
    # stray "Text" number format applied to the cell.
    $cell.Style = "Normal"
}
